$wb = $excel.ActiveWorkbook

# --- Add "Chapter 4" sheet (empty) after Chapter 3 ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws5.Name = "Chapter 4"

# --- Add "The rest" sheet after Chapter 4 ---
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws6 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet2)
$ws6.Name = "The rest"

# --- Column width for "The rest" ---
$ws6.Columns.Item(1).ColumnWidth = 63.42578125

# --- Fill "The rest" content (order matters: matches shared-string build order) ---
$ws6.Range("A16").Value = " Paka ichowo Bura nen Adech III, Siro 2"
$ws6.Range("A15").Value = "  To chowo Bura gikenyo; to wok ndirino to Majanga bedo ja tel ma chowo Bura"
$ws6.Range("A14").Value = "  Bura bedo i tele man yo Nyawiyoga"
$ws6.Range("A13").Value = "Yado chon onwaŋo"
$ws6.Range("A12").Value = "  Kumo okel gigipiny me totero jo yo lul Tewo (Nyakiriga) towacho rijo ni me atele ma kabedo manyien pa Bura "
$ws6.Range("A11").Value = "  To munyo onwaŋ, ti penjo kumago owok iye, to kwero wacho rijo kwanyo woko ni kajo kelo rigo koŋo kod wot gwendi"
$ws6.Range("A10").Value = "I wacho ni nitye ndir ma Majanga orwenyo makinen tieko ndelo maromo abiriyo"
$ws6.Range("A9").Value = "  Ri ameno to ji gye chako lworo Majanga gi miyogo dwoŋ gi winjo wach pere"
$ws6.Range("A8").Value = "  I lweny bende yado kowacho ni kidh win kumanyo, ko kidh thenge no tituro lweny, to kowacho ni kononi kada wikidho i bino turo win, apaka bende bedo"
$ws6.Range("A7").Value = "  Odoko ruman owacho ri ji ni nitye ogwaŋ Kwach neko nyako moro ka moro, ka ni kidho kenyo onwaŋ nyako no otieko kir tho"
$ws6.Range("A6").Value = "  Wadi jo ko yeyo; to rigiraura madwoŋ ka ni kidho poyo, Nyielo to neko ŋato"
$ws6.Range("A5").Value = "  Obedo ndir moro kuma jo kidho dwar, Majanga owacho ni nyawoti gin achiel Nyielo ya neko"
$ws6.Range("A4").Value = "To wok chon chango Majanga oneno paka nitye kod gi moro iwiye pa jwok, ma kowacho ni gi moro ya timere to gino timere atima ameno"
$ws6.Range("A3").Value = "  To rupir chango baa mere otii aka ja chandi, omin baa mere mi lwoŋo ni Akure amunywomo rigo; kareno kutho Majanga odongo pere kakwayo rigo dhoki aka omito go swa pa nyath pere won"
$ws6.Range("A2").Value = "Majanga chango baa mere i lwoŋo ni Kinara ja Nyapolo Ogule"
$ws6.Range("A18").Value = "  Kendo chango obedo ja lweny mamisen swa, ama bin telo kir ji yo lweny"
$ws6.Range("A19").Value = "  Kir Akisili gye (nen Adech III, Siro 8), to bedo jatel mere"
$ws6.Range("A21").Value = "  To kada ameno otemo wacho ri jii ni joyikere ri lweny no; to rupir jo nicha chango jo lwenyo gi mundu amumiyo otur jo ma piyo piyo"
$ws6.Range("A22").Value = " Rumachien pa bino pa Kakunguru, Majanga obedo hongo manok to tho i oro chiegin 1905"
$ws6.Range("A23").Value = "SIRO 5"
$ws6.Range("A17").Value = "Majanga to limo dwoŋ madit ri kwom gigipiny me gye"
$ws6.Range("A20").Value = " Ndiri ma Kakunguru donjo i piny me gi turo, onwa?o Majanga onyo obedo jadwoŋ m'oti"

# --- Selections on existing sheets (content / view state) ---
$wsAck = $wb.Worksheets.Item("Acknowledgement and Dedication")
$wsAck.Range("B34").Select()

$wsCh1 = $wb.Worksheets.Item("Chapter 1")
$wsCh1.Range("A1:B1").Select()

$wsCh2 = $wb.Worksheets.Item("Chapter 2")
$wsCh2.Range("B61").Select()

$wsCh3 = $wb.Worksheets.Item("Chapter 3")
$wsCh3.Range("B22").Select()

# --- Selections on new sheets ---
$ws5.Range("I20").Select()

$ws6.Range("F14").Select()

# --- Final active sheet = "The rest" (matches activeTab=5) ---
$ws6.Activate()

